# In the Mix - Two Day Schedule: rename two session titles to "Data Science"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Day 1, 11:15 AM session: "Computational, Part 1" -> "Data Science, Part 1"
$ws.Range("C11").Value = "Data Science, Part 1"

# Day 2, 3:45 PM session: "Machine Learning, Part 2" -> "Data Science, Part 2"
$ws.Range("C29").Value = "Data Science, Part 2"
